# Apply scheduled-runner price/profit updates to each Sheet's leve profitability table.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H41").Value = 2182.4666
$ws.Range("I41").Value = 1534.1666
$ws.Range("K41").Value = 1534.1666
$ws.Range("M41").Value = -1094.1666
$ws.Range("H137").Value = 1037.7858
$ws.Range("I137").Value = 937.1818
$ws.Range("K137").Value = 2811.5454
$ws.Range("M137").Value = -261.5454
$ws.Range("H141").Value = 7404
$ws.Range("I141").Value = 7577.091
$ws.Range("K141").Value = 22731.273
$ws.Range("M141").Value = -17551.273

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 4014.275
$ws.Range("I61").Value = 2971.3333
$ws.Range("K61").Value = 2971.3333
$ws.Range("M61").Value = -2759.3333
$ws.Range("H74").Value = 27102.027
$ws.Range("I74").Value = 1224.7778
$ws.Range("K74").Value = 1224.7778
$ws.Range("M74").Value = -350.7778000000001
$ws.Range("H77").Value = 27102.027
$ws.Range("I77").Value = 1224.7778
$ws.Range("K77").Value = 6123.889
$ws.Range("M77").Value = -1755.889
$ws.Range("H136").Value = 4014.275
$ws.Range("I136").Value = 2971.3333
$ws.Range("K136").Value = 8913.999899999999
$ws.Range("M136").Value = -6363.999899999999

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H105").Value = 62502676
$ws.Range("I105").Value = 100001960
$ws.Range("J105").Value = 3866.6667
$ws.Range("K105").Value = 100001960
$ws.Range("L105").Value = 3866.6667
$ws.Range("M105").Value = -100000213
$ws.Range("N105").Value = -7360.6667

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 3179.3333
$ws.Range("I16").Value = 3461.5557
$ws.Range("K16").Value = 3461.5557
$ws.Range("M16").Value = -3174.5557
$ws.Range("H31").Value = 1499.7646
$ws.Range("I31").Value = 1482.8422
$ws.Range("K31").Value = 1482.8422
$ws.Range("M31").Value = -1187.8422
$ws.Range("H33").Value = 1513.3334
$ws.Range("I33").Value = 1513.3334
$ws.Range("K33").Value = 1513.3334
$ws.Range("M33").Value = -1134.3334
$ws.Range("H34").Value = 1499.7646
$ws.Range("I34").Value = 1482.8422
$ws.Range("K34").Value = 1482.8422
$ws.Range("M34").Value = -1280.8422
$ws.Range("H58").Value = 43862904
$ws.Range("I58").Value = 33336564
$ws.Range("K58").Value = 33336564
$ws.Range("M58").Value = -33336361
$ws.Range("H113").Value = 3179.3333
$ws.Range("I113").Value = 3461.5557
$ws.Range("K113").Value = 3461.5557
$ws.Range("M113").Value = -1291.5557
$ws.Range("H132").Value = 6182.4165
$ws.Range("I132").Value = 4354.5557
$ws.Range("K132").Value = 13063.6671
$ws.Range("M132").Value = -10533.6671
$ws.Range("H134").Value = 3523.75
$ws.Range("I134").Value = 3523.75
$ws.Range("K134").Value = 10571.25
$ws.Range("M134").Value = -8036.25
$ws.Range("H136").Value = 43862904
$ws.Range("I136").Value = 33336564
$ws.Range("K136").Value = 100009692
$ws.Range("M136").Value = -100007142

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H33").Value = 426
$ws.Range("J33").Value = 1035
$ws.Range("L33").Value = 6210
$ws.Range("N33").Value = -6776
$ws.Range("H34").Value = 966.26666
$ws.Range("I34").Value = 309.7
$ws.Range("J34").Value = 2279.4
$ws.Range("K34").Value = 929.0999999999999
$ws.Range("L34").Value = 6838.200000000001
$ws.Range("M34").Value = -845.0999999999999
$ws.Range("N34").Value = -7006.200000000001
$ws.Range("H39").Value = 8191.1
$ws.Range("I39").Value = 7500
$ws.Range("J39").Value = 8267.888999999999
$ws.Range("K39").Value = 22500
$ws.Range("L39").Value = 24803.667
$ws.Range("M39").Value = -22206
$ws.Range("N39").Value = -25391.667
$ws.Range("H55").Value = 13447.333
$ws.Range("I55").Value = 2500
$ws.Range("J55").Value = 14442.546
$ws.Range("K55").Value = 7500
$ws.Range("L55").Value = 43327.638
$ws.Range("M55").Value = -7323
$ws.Range("N55").Value = -43681.638
$ws.Range("H107").Value = 3397.7273
$ws.Range("J107").Value = 3660.5
$ws.Range("L107").Value = 10981.5
$ws.Range("N107").Value = -14821.5
$ws.Range("H117").Value = 433948.06
$ws.Range("I117").Value = 1941.8334
$ws.Range("J117").Value = 693151.8
$ws.Range("K117").Value = 5825.5002
$ws.Range("L117").Value = 2079455.4
$ws.Range("M117").Value = -2383.5002
$ws.Range("N117").Value = -2086339.4

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H97").Value = 30226.594
$ws.Range("I97").Value = 49388.895
$ws.Range("J97").Value = 2220.1538
$ws.Range("K97").Value = 49388.895
$ws.Range("L97").Value = 2220.1538
$ws.Range("M97").Value = -48892.895
$ws.Range("N97").Value = -3212.1538
$ws.Range("H102").Value = 54048.41
$ws.Range("I102").Value = 103538.2
$ws.Range("J102").Value = 12806.917
$ws.Range("K102").Value = 103538.2
$ws.Range("L102").Value = 12806.917
$ws.Range("M102").Value = -101916.2
$ws.Range("N102").Value = -16050.917
$ws.Range("H107").Value = 1164.1765
$ws.Range("I107").Value = 1026.909
$ws.Range("K107").Value = 1026.909
$ws.Range("M107").Value = 893.0909999999999
$ws.Range("H113").Value = 2976.9092
$ws.Range("J113").Value = 3097.1667
$ws.Range("L113").Value = 3097.1667
$ws.Range("N113").Value = -7437.1667
$ws.Range("H132").Value = 8890.267
$ws.Range("I132").Value = 7336.1
$ws.Range("K132").Value = 22008.3
$ws.Range("M132").Value = -19478.3

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H40").Value = 4987.533
$ws.Range("I40").Value = 4835.231
$ws.Range("K40").Value = 4835.231
$ws.Range("M40").Value = -4699.231
$ws.Range("H46").Value = 2797.9812
$ws.Range("I46").Value = 2199
$ws.Range("J46").Value = 2809.5
$ws.Range("K46").Value = 2199
$ws.Range("L46").Value = 2809.5
$ws.Range("M46").Value = -2011
$ws.Range("N46").Value = -3185.5
$ws.Range("H82").Value = 1527
$ws.Range("I82").Value = 1577.4
$ws.Range("J82").Value = 1383
$ws.Range("K82").Value = 1577.4
$ws.Range("L82").Value = 1383
$ws.Range("M82").Value = -1216.4
$ws.Range("N82").Value = -2105
$ws.Range("H85").Value = 1527
$ws.Range("I85").Value = 1577.4
$ws.Range("J85").Value = 1383
$ws.Range("K85").Value = 1577.4
$ws.Range("L85").Value = 1383
$ws.Range("M85").Value = -329.4000000000001
$ws.Range("N85").Value = -3879
$ws.Range("H100").Value = 3069.625
$ws.Range("I100").Value = 2484.5
$ws.Range("J100").Value = 4044.8333
$ws.Range("K100").Value = 2484.5
$ws.Range("L100").Value = 4044.8333
$ws.Range("M100").Value = -1943.5
$ws.Range("N100").Value = -5126.8333
$ws.Range("H132").Value = 3211.6667
$ws.Range("I132").Value = 2678
$ws.Range("K132").Value = 8034
$ws.Range("M132").Value = -5504
$ws.Range("H136").Value = 58964.688
$ws.Range("I136").Value = 2895.6667
$ws.Range("J136").Value = 900000
$ws.Range("K136").Value = 8687.000100000001
$ws.Range("L136").Value = 2700000
$ws.Range("M136").Value = -6137.000100000001
$ws.Range("N136").Value = -2705100

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H107").Value = 338.45456
$ws.Range("I107").Value = 302.66666
$ws.Range("K107").Value = 907.9999799999999
$ws.Range("M107").Value = 1012.00002
$ws.Range("H126").Value = 3267.0312
$ws.Range("I126").Value = 3366
$ws.Range("J126").Value = 2970.125
$ws.Range("K126").Value = 10098
$ws.Range("L126").Value = 8910.375
$ws.Range("M126").Value = -7628
$ws.Range("N126").Value = -13850.375
